$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns retain text formatting so values are not
# auto-converted to numbers by Excel (the source data stores these as text).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '29.096.67'
$ws.Range('E2').Value = '  +0.25%  '
$ws.Range('D3').Value = '1.837.02'
$ws.Range('E3').Value = '  +0.50%  '
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.27%  '
$ws.Range('D5').Value = '244.78'
$ws.Range('E5').Value = '  +1.69%  '
$ws.Range('D6').Value = '0.6366'
$ws.Range('E6').Value = '  +2.55%  '
$ws.Range('E7').Value = '  +0.15%  '
$ws.Range('D8').Value = '0.07568'
$ws.Range('E8').Value = '  +1.04%  '
$ws.Range('D9').Value = '0.2953'
$ws.Range('E9').Value = '  +1.36%  '
$ws.Range('D10').Value = '22.93'
$ws.Range('E10').Value = '  +1.33%  '
$ws.Range('D11').Value = '0.07752'
$ws.Range('E11').Value = '  +1.82%  '
$ws.Range('D12').Value = '1.841.44'
$ws.Range('E12').Value = '  +0.82%  '
$ws.Range('D13').Value = '5.012'
$ws.Range('E13').Value = '  +1.37%  '
$ws.Range('D14').Value = '0.6734'
$ws.Range('E14').Value = '  +1.82%  '
$ws.Range('D15').Value = '83.46'
$ws.Range('E15').Value = '  +1.83%  '
$ws.Range('D16').Value = '0.000009561'
$ws.Range('E16').Value = '  +5.75%  '
$ws.Range('D17').Value = '6.125'
$ws.Range('E17').Value = '  +2.86%  '
$ws.Range('D18').Value = '29.131.67'
$ws.Range('D19').Value = '12.63'
$ws.Range('E19').Value = '  +2.61%  '
$ws.Range('D20').Value = '227.23'
$ws.Range('E20').Value = '  +1.38%  '
$ws.Range('E21').Value = '  +0.09%  '
$ws.Range('D22').Value = '7.226'
$ws.Range('E22').Value = '  +0.91%  '
$ws.Range('D23').Value = '1.003'
$ws.Range('E23').Value = '  +0.23%  '
$ws.Range('D24').Value = '160.87'
$ws.Range('E24').Value = '  +0.92%  '
$ws.Range('D25').Value = '0.1403'
$ws.Range('E25').Value = '  +3.75%  '
$ws.Range('D26').Value = '8.560'
$ws.Range('E26').Value = '  +2.21%  '
$ws.Range('D27').Value = '17.99'
$ws.Range('E27').Value = '  +1.19%  '
$ws.Range('D28').Value = '1.505'
$ws.Range('E28').Value = '  +0.86%  '
$ws.Range('D29').Value = '4.132'
$ws.Range('E29').Value = '  +2.34%  '
$ws.Range('D30').Value = '4.087'
$ws.Range('E30').Value = '  +1.52%  '
$ws.Range('D31').Value = '1.206'
$ws.Range('E31').Value = '  +0.21%  '
$ws.Range('E32').Value = '  +3.68%  '
$ws.Range('D33').Value = '1.868'
$ws.Range('E33').Value = '  +2.16%  '
$ws.Range('D34').Value = '0.7495'
$ws.Range('E34').Value = '  +2.57%  '
$ws.Range('D35').Value = '1.144'
$ws.Range('E35').Value = '  -0.27%  '
$ws.Range('D36').Value = '2.662'
$ws.Range('E36').Value = '  +0.59%  '
$ws.Range('D37').Value = '1.244.78'
$ws.Range('E37').Value = '  -2.38%  '
$ws.Range('D38').Value = '2.765'
$ws.Range('E38').Value = '  +0.59%  '
$ws.Range('D39').Value = '0.01791'
$ws.Range('E39').Value = '  +0.76%  '
$ws.Range('D40').Value = '6.630'
$ws.Range('E40').Value = '  +5.13%  '
$ws.Range('D41').Value = '0.9090'
$ws.Range('E41').Value = '  +2.06%  '
$ws.Range('E42').Value = '  +0.04%  '
$ws.Range('D43').Value = '102.34'
$ws.Range('E43').Value = '  +0.79%  '
$ws.Range('D44').Value = '1.990.75'
$ws.Range('E44').Value = '  +0.79%  '
$ws.Range('B45').Value = 'BabyDogeCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D45').Value = '0.00000000126'
$ws.Range('E45').Value = '  +6.04%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = '65.31'
$ws.Range('E46').Value = '  +3.13%  '
$ws.Range('D47').Value = '0.5117'
$ws.Range('E47').Value = '  +0.06%  '
$ws.Range('D48').Value = '0.4105'
$ws.Range('E48').Value = '  +3.84%  '
$ws.Range('D49').Value = '9.139'
$ws.Range('E49').Value = '  +3.76%  '
$ws.Range('D50').Value = '6.793'
$ws.Range('E50').Value = '  +2.29%  '
$ws.Range('D51').Value = '1.653'
$ws.Range('E51').Value = '  -2.08%  '
